# Auto-generated: apply the cryptos.xlsx price/volume update described by the diff.
# D column (Price) holds number-looking text (e.g. "1.001", "0.00000000117",
# "28.843.90") that must stay literal TEXT -- plain `.Value =` assignment lets
# Excel's COM layer auto-coerce these into floating point numbers, which both
# changes their type and silently rewrites their displayed digits (trailing
# zeros dropped, thousands-separator-looking values collapsed, exponential
# notation for tiny numbers, etc). To avoid that, D-column cells are forced to
# Text format before the value is written, then restored to the default
# "Normal" style afterwards so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '28.843.90'
$ws.Range("E2").Value = '  -1.64%  '
Set-TextValue "D3" '1.830.15'
$ws.Range("E3").Value = '  -1.78%  '
Set-TextValue "D4" '0.9999'
$ws.Range("E4").Value = '  -0.05%  '
Set-TextValue "D5" '244.47'
$ws.Range("E5").Value = '  +0.59%  '
Set-TextValue "D6" '0.6925'
$ws.Range("E6").Value = '  -1.17%  '
Set-TextValue "D7" '1.001'
$ws.Range("E7").Value = '  +0.01%  '
Set-TextValue "D8" '0.07668'
$ws.Range("E8").Value = '  -3.04%  '
Set-TextValue "D9" '0.3041'
$ws.Range("E9").Value = '  -2.63%  '
Set-TextValue "D10" '23.18'
$ws.Range("E10").Value = '  -4.77%  '
Set-TextValue "D11" '0.07788'
$ws.Range("E11").Value = '  -0.40%  '
Set-TextValue "D12" '93.22'
$ws.Range("E12").Value = '  +0.95%  '
Set-TextValue "D13" '1.833.42'
$ws.Range("E13").Value = '  -1.84%  '
Set-TextValue "D14" '5.082'
$ws.Range("E14").Value = '  -1.18%  '
Set-TextValue "D15" '0.6775'
$ws.Range("E15").Value = '  -2.88%  '
Set-TextValue "D16" '6.436'
$ws.Range("E16").Value = '  -1.62%  '
Set-TextValue "D17" '0.000008234'
$ws.Range("E17").Value = '  -3.60%  '
Set-TextValue "D18" '28.853.50'
$ws.Range("E18").Value = '  -1.69%  '
Set-TextValue "D19" '241.95'
$ws.Range("E19").Value = '  -2.82%  '
$ws.Range("E20").Value = '  -2.36%  '
Set-TextValue "D21" '12.66'
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("E22").Value = '  +0.06%  '
Set-TextValue "D23" '7.426'
$ws.Range("E23").Value = '  -2.12%  '
Set-TextValue "D24" '1.001'
$ws.Range("E24").Value = '  +0.02%  '
Set-TextValue "D25" '0.1486'
$ws.Range("E25").Value = '  -3.24%  '
Set-TextValue "D26" '159.99'
$ws.Range("E26").Value = '  -0.47%  '
Set-TextValue "D27" '8.743'
$ws.Range("E27").Value = '  -2.46%  '
Set-TextValue "D28" '18.21'
$ws.Range("E28").Value = '  -2.79%  '
Set-TextValue "D29" '1.541'
$ws.Range("E29").Value = '  -2.63%  '
Set-TextValue "D30" '4.213'
$ws.Range("E30").Value = '  -2.01%  '
Set-TextValue "D31" '4.150'
$ws.Range("E31").Value = '  -2.17%  '
Set-TextValue "D32" '1.182'
$ws.Range("E32").Value = '  -2.12%  '
Set-TextValue "D33" '0.05093'
$ws.Range("E33").Value = '  -2.73%  '
Set-TextValue "D34" '0.7717'
$ws.Range("E34").Value = '  +1.95%  '
Set-TextValue "D35" '1.857'
$ws.Range("E35").Value = '  -1.53%  '
Set-TextValue "D36" '1.138'
$ws.Range("E36").Value = '  -3.32%  '
Set-TextValue "D37" '2.695'
$ws.Range("E37").Value = '  -0.15%  '
Set-TextValue "D38" '0.01846'
$ws.Range("E38").Value = '  -1.03%  '
Set-TextValue "D39" '1.239.93'
$ws.Range("E39").Value = '  -2.79%  '
Set-TextValue "D40" '2.700'
$ws.Range("E40").Value = '  -1.66%  '
Set-TextValue "D41" '0.9514'
$ws.Range("E41").Value = '  +5.49%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D42" '5.979'
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D43" '107.49'
$ws.Range("E43").Value = '  -2.49%  '
Set-TextValue "D44" '1.001'
$ws.Range("E44").Value = '  +0.06%  '
Set-TextValue "D45" '9.616'
$ws.Range("E45").Value = '  -0.02%  '
Set-TextValue "D46" '1.976.24'
$ws.Range("E46").Value = '  -2.27%  '
Set-TextValue "D47" '0.5155'
$ws.Range("E47").Value = '  -0.36%  '
Set-TextValue "D48" '63.84'
$ws.Range("E48").Value = '  -8.92%  '
Set-TextValue "D49" '1.738'
$ws.Range("E49").Value = '  -2.85%  '
Set-TextValue "D50" '0.00000000117'
$ws.Range("E50").Value = '  -6.37%  '
Set-TextValue "D51" '6.908'
$ws.Range("E51").Value = '  -1.55%  '
